# Implementacion de Dublin Core en las tablas de metadatos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Referentes")

# Update header row (row 1) of "Referentes" sheet to use Dublin Core field names
$ws.Range("B1").Value = "dc.title"
$ws.Range("C1").Value = "dc.date"
$ws.Range("D1").Value = "dc.publisher"
$ws.Range("A1").Value = "dc.identifier"

# Update the active selection on this sheet
$ws.Range("B10").Select()
